$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.86326320006514
$ws.Range("C2").Value = 14.03311813798708
$ws.Range("D2").Value = 4.67553686572749
$ws.Range("E2").Value = 16.49936937945428
$ws.Range("F2").Value = 30.04103356718291
$ws.Range("I2").Value = 21.78920218704911
$ws.Range("N2").Value = 16.63493714321157
$ws.Range("B3").Value = 16.09950608411895
$ws.Range("C3").Value = 13.20107851691508
$ws.Range("D3").Value = 4.705322109001955
$ws.Range("E3").Value = 15.55653740869772
$ws.Range("F3").Value = 29.53962476997393
$ws.Range("I3").Value = 21.71450993813013
$ws.Range("N3").Value = 16.71323766419159
$ws.Range("B4").Value = 15.6165365790916
$ws.Range("C4").Value = 12.66614278142704
$ws.Range("D4").Value = 4.724737553455235
$ws.Range("E4").Value = 14.9537146056006
$ws.Range("F4").Value = 29.24008450748608
$ws.Range("I4").Value = 21.67630736357052
$ws.Range("N4").Value = 16.7632661077178
$ws.Range("B5").Value = 15.4165069466694
$ws.Range("C5").Value = 12.44229443434551
$ws.Range("D5").Value = 4.732931184931721
$ws.Range("E5").Value = 14.70231493588561
$ws.Range("F5").Value = 29.12027612849793
$ws.Range("I5").Value = 21.66266416984451
$ws.Range("N5").Value = 16.78414567939432
$ws.Range("B6").Value = 15.38310782455767
$ws.Range("C6").Value = 12.40477725121553
$ws.Range("D6").Value = 4.734308694129534
$ws.Range("E6").Value = 14.66023221925074
$ws.Range("F6").Value = 29.10052298770863
$ws.Range("I6").Value = 21.66051488219825
$ws.Range("N6").Value = 16.78764252021187
$ws.Range("B7").Value = 15.61385148993465
$ws.Range("C7").Value = 12.66314731696181
$ws.Range("D7").Value = 4.724846917373409
$ws.Range("E7").Value = 14.95034700444685
$ws.Range("F7").Value = 29.23845938126506
$ws.Range("I7").Value = 21.67611557741753
$ws.Range("N7").Value = 16.76354570006199
$ws.Range("B8").Value = 16.60300276344202
$ws.Range("C8").Value = 13.75132760860096
$ws.Range("D8").Value = 4.685571788110286
$ws.Range("E8").Value = 16.17939195783057
$ws.Range("F8").Value = 29.86651887182369
$ws.Range("I8").Value = 21.76185885306386
$ws.Range("N8").Value = 16.66153121405278
$ws.Range("B9").Value = 18.41988254978029
$ws.Range("C9").Value = 15.68778390903621
$ws.Range("D9").Value = 4.617579445633492
$ws.Range("E9").Value = 18.48260987542155
$ws.Range("F9").Value = 31.15683926971352
$ws.Range("I9").Value = 21.99069407608192
$ws.Range("N9").Value = 16.47688076695082
$ws.Range("B10").Value = 19.6668466125248
$ws.Range("C10").Value = 16.9834811505773
$ws.Range("D10").Value = 4.57324683375428
$ws.Range("E10").Value = 20.14762134080832
$ws.Range("F10").Value = 32.13060688052058
$ws.Range("I10").Value = 22.19553433515441
$ws.Range("N10").Value = 16.35048352457731
$ws.Range("B11").Value = 20.21292546750902
$ws.Range("C11").Value = 17.54442871298576
$ws.Range("D11").Value = 4.554327269562337
$ws.Range("E11").Value = 20.8636324526438
$ws.Range("F11").Value = 32.57707905736765
$ws.Range("I11").Value = 22.29657491870482
$ws.Range("N11").Value = 16.29496772538405
$ws.Range("B12").Value = 20.41652675317806
$ws.Range("C12").Value = 17.75269737804442
$ws.Range("D12").Value = 4.547344847858403
$ws.Range("E12").Value = 21.12885469364625
$ws.Range("F12").Value = 32.746477391517
$ws.Range("I12").Value = 22.33595195382016
$ws.Range("N12").Value = 16.27422853150778
$ws.Range("B13").Value = 20.37282133760514
$ws.Range("C13").Value = 17.70802830472552
$ws.Range("D13").Value = 4.548840499488057
$ws.Range("E13").Value = 21.07199681444409
$ws.Range("F13").Value = 32.70998251709078
$ws.Range("I13").Value = 22.32742206977603
$ws.Range("N13").Value = 16.27868250921821
$ws.Range("B14").Value = 20.22974047279995
$ws.Range("C14").Value = 17.56164654247367
$ws.Range("D14").Value = 4.553749158419958
$ws.Range("E14").Value = 20.88557077964896
$ws.Range("F14").Value = 32.59100994671133
$ws.Range("I14").Value = 22.29979223527027
$ws.Range("N14").Value = 16.29325582873647
$ws.Range("B15").Value = 20.14168044635809
$ws.Range("C15").Value = 17.47144162304754
$ws.Range("D15").Value = 4.55677963612174
$ws.Range("E15").Value = 20.77061022358257
$ws.Range("F15").Value = 32.51817357174709
$ws.Range("I15").Value = 22.28301290865018
$ws.Range("N15").Value = 16.30221926869253
$ws.Range("B16").Value = 19.6307203331286
$ws.Range("C16").Value = 16.94624387798889
$ws.Range("D16").Value = 4.574508595819851
$ws.Range("E16").Value = 20.09999839105503
$ws.Range("F16").Value = 32.10148429507814
$ws.Range("I16").Value = 22.18908800856398
$ws.Range("N16").Value = 16.35415136799237
$ws.Range("B17").Value = 19.31173138345702
$ws.Range("C17").Value = 16.6167153150177
$ws.Range("D17").Value = 4.585706082447154
$ws.Range("E17").Value = 19.67802134919618
$ws.Range("F17").Value = 31.84662989864033
$ws.Range("I17").Value = 22.13347098353924
$ws.Range("N17").Value = 16.38651668687908
$ws.Range("B18").Value = 19.12627153843833
$ws.Range("C18").Value = 16.42450161079787
$ws.Range("D18").Value = 4.592263842929481
$ws.Range("E18").Value = 19.43141118104568
$ws.Range("F18").Value = 31.70038403248875
$ws.Range("I18").Value = 22.10222220004287
$ws.Range("N18").Value = 16.40531910759829
$ws.Range("B19").Value = 19.06314168284407
$ws.Range("C19").Value = 16.35896341945634
$ws.Range("D19").Value = 4.59450425055316
$ws.Range("E19").Value = 19.34724216780952
$ws.Range("F19").Value = 31.65093142035202
$ws.Range("I19").Value = 22.09176951546285
$ws.Range("N19").Value = 16.41171741009099
$ws.Range("B20").Value = 19.34589493503634
$ws.Range("C20").Value = 16.65207176417217
$ws.Range("D20").Value = 4.584501936216114
$ws.Range("E20").Value = 19.72334489026741
$ws.Range("F20").Value = 31.87372571133859
$ws.Range("I20").Value = 22.13931495439393
$ws.Range("N20").Value = 16.38305202843518
$ws.Range("B21").Value = 20.27185430467744
$ws.Range("C21").Value = 17.6047554181689
$ws.Range("D21").Value = 4.552302405222759
$ws.Range("E21").Value = 20.94048886148515
$ws.Range("F21").Value = 32.62594750916708
$ws.Range("I21").Value = 22.30787765259198
$ws.Range("N21").Value = 16.28896761276249
$ws.Range("B22").Value = 20.85839602782384
$ws.Range("C22").Value = 18.20328529036729
$ws.Range("D22").Value = 4.532320332899815
$ws.Range("E22").Value = 21.70150014141525
$ws.Range("F22").Value = 33.11940869294305
$ws.Range("I22").Value = 22.42453451260911
$ws.Range("N22").Value = 16.22912929585626
$ws.Range("B23").Value = 20.54709322075044
$ws.Range("C23").Value = 17.88602087264057
$ws.Range("D23").Value = 4.542887053555599
$ws.Range("E23").Value = 21.29847450150497
$ws.Range("F23").Value = 32.85592695077821
$ws.Range("I23").Value = 22.36168411516848
$ws.Range("N23").Value = 16.26091560062929
$ws.Range("B24").Value = 19.33045602631613
$ws.Range("C24").Value = 16.63609571632843
$ws.Range("D24").Value = 4.58504595660068
$ws.Range("E24").Value = 19.70286661040437
$ws.Range("F24").Value = 31.86147482979894
$ws.Range("I24").Value = 22.13667063217627
$ws.Range("N24").Value = 16.38461779248437
$ws.Range("B25").Value = 17.94299181732663
$ws.Range("C25").Value = 15.18584752800484
$ws.Range("D25").Value = 4.634995031346753
$ws.Range("E25").Value = 17.83247183317552
$ws.Range("F25").Value = 30.80252218164245
$ws.Range("I25").Value = 21.92230873480763
$ws.Range("N25").Value = 16.52519767638384
